# Auto-generated edit script applying scraped price-refresh diff
# to Titan_Profits workbook (8 item-category sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Cells.Item(10, 10).Value = 50000  # J10: 0 -> 50000
$ws.Cells.Item(10, 12).Value = 50000  # L10: 0 -> 50000
$ws.Cells.Item(10, 14).Value = -50586  # N10: None -> -50586
# Row 15
$ws.Cells.Item(15, 8).Value = 157608.4  # H15: 102478.47 -> 157608.4
$ws.Cells.Item(15, 9).Value = 157608.4  # I15: 102478.47 -> 157608.4
$ws.Cells.Item(15, 11).Value = 472825.2  # K15: 307435.41 -> 472825.2
$ws.Cells.Item(15, 13).Value = -472656.2  # M15: -307266.41 -> -472656.2
# Row 95
$ws.Cells.Item(95, 8).Value = 30000  # H95: 0 -> 30000
$ws.Cells.Item(95, 10).Value = 30000  # J95: 0 -> 30000
$ws.Cells.Item(95, 12).Value = 30000  # L95: 0 -> 30000
$ws.Cells.Item(95, 14).Value = -35492  # N95: None -> -35492
# Row 132
$ws.Cells.Item(132, 8).Value = 33688.773  # H132: 40875.848 -> 33688.773
$ws.Cells.Item(132, 9).Value = 33688.773  # I132: 42390.88 -> 33688.773
$ws.Cells.Item(132, 10).Value = 0  # J132: 3000 -> 0
$ws.Cells.Item(132, 11).Value = 101066.319  # K132: 127172.64 -> 101066.319
$ws.Cells.Item(132, 12).Value = 0  # L132: 9000 -> 0
$ws.Cells.Item(132, 13).Value = -98536.319  # M132: -124642.64 -> -98536.319
$ws.Cells.Item(132, 14).ClearContents()  # N132: -14060 -> (cleared)

$ws = $wb.Worksheets.Item("ARM")
# Row 25
$ws.Cells.Item(25, 8).Value = 0  # H25: 2000 -> 0
$ws.Cells.Item(25, 9).Value = 0  # I25: 2000 -> 0
$ws.Cells.Item(25, 11).Value = 0  # K25: 2000 -> 0
$ws.Cells.Item(25, 13).ClearContents()  # M25: -1598 -> (cleared)
# Row 74
$ws.Cells.Item(74, 8).Value = 5718.7144  # H74: 6155.5386 -> 5718.7144
$ws.Cells.Item(74, 9).Value = 1067.625  # I74: 1154 -> 1067.625
$ws.Cells.Item(74, 10).Value = 11920.167  # J74: 11157.077 -> 11920.167
$ws.Cells.Item(74, 11).Value = 1067.625  # K74: 1154 -> 1067.625
$ws.Cells.Item(74, 12).Value = 11920.167  # L74: 11157.077 -> 11920.167
$ws.Cells.Item(74, 13).Value = -193.625  # M74: -280 -> -193.625
$ws.Cells.Item(74, 14).Value = -13668.167  # N74: -12905.077 -> -13668.167
# Row 77
$ws.Cells.Item(77, 8).Value = 5718.7144  # H77: 6155.5386 -> 5718.7144
$ws.Cells.Item(77, 9).Value = 1067.625  # I77: 1154 -> 1067.625
$ws.Cells.Item(77, 10).Value = 11920.167  # J77: 11157.077 -> 11920.167
$ws.Cells.Item(77, 11).Value = 5338.125  # K77: 5770 -> 5338.125
$ws.Cells.Item(77, 12).Value = 59600.835  # L77: 55785.38499999999 -> 59600.835
$ws.Cells.Item(77, 13).Value = -970.125  # M77: -1402 -> -970.125
$ws.Cells.Item(77, 14).Value = -68336.83499999999  # N77: -64521.38499999999 -> -68336.83499999999
# Row 110
$ws.Cells.Item(110, 8).Value = 1035.84  # H110: 1392.4445 -> 1035.84
$ws.Cells.Item(110, 9).Value = 699.15  # I110: 1060.2727 -> 699.15
$ws.Cells.Item(110, 10).Value = 2382.6  # J110: 1914.4286 -> 2382.6
$ws.Cells.Item(110, 11).Value = 699.15  # K110: 1060.2727 -> 699.15
$ws.Cells.Item(110, 12).Value = 2382.6  # L110: 1914.4286 -> 2382.6
$ws.Cells.Item(110, 13).Value = 1345.85  # M110: 984.7273 -> 1345.85
$ws.Cells.Item(110, 14).Value = -6472.6  # N110: -6004.4286 -> -6472.6
# Row 122
$ws.Cells.Item(122, 8).Value = 2781.0322  # H122: 2866.6667 -> 2781.0322
$ws.Cells.Item(122, 9).Value = 2052.6956  # I122: 2136.3635 -> 2052.6956
$ws.Cells.Item(122, 11).Value = 6158.0868  # K122: 6409.0905 -> 6158.0868
$ws.Cells.Item(122, 13).Value = -3708.0868  # M122: -3959.0905 -> -3708.0868

$ws = $wb.Worksheets.Item("BSM")
# Row 29
$ws.Cells.Item(29, 8).Value = 750  # H29: 1505.3334 -> 750
$ws.Cells.Item(29, 9).Value = 750  # I29: 1505.3334 -> 750
$ws.Cells.Item(29, 11).Value = 750  # K29: 1505.3334 -> 750
$ws.Cells.Item(29, 13).Value = -461  # M29: -1216.3334 -> -461
# Row 134
$ws.Cells.Item(134, 8).Value = 3935.7334  # H134: 4337.6 -> 3935.7334
$ws.Cells.Item(134, 9).Value = 2070.9  # I134: 2243.5 -> 2070.9
$ws.Cells.Item(134, 10).Value = 7665.4  # J134: 7478.75 -> 7665.4
$ws.Cells.Item(134, 11).Value = 6212.700000000001  # K134: 6730.5 -> 6212.700000000001
$ws.Cells.Item(134, 12).Value = 22996.2  # L134: 22436.25 -> 22996.2
$ws.Cells.Item(134, 13).Value = -3677.700000000001  # M134: -4195.5 -> -3677.700000000001
$ws.Cells.Item(134, 14).Value = -28066.2  # N134: -27506.25 -> -28066.2

$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Cells.Item(86, 8).Value = 55558852  # H86: 50002000 -> 55558852
$ws.Cells.Item(86, 9).Value = 250001000  # I86: 71429710 -> 250001000
$ws.Cells.Item(86, 10).Value = 3953.5715  # J86: 3999.6667 -> 3953.5715
$ws.Cells.Item(86, 11).Value = 250001000  # K86: 71429710 -> 250001000
$ws.Cells.Item(86, 12).Value = 3953.5715  # L86: 3999.6667 -> 3953.5715
$ws.Cells.Item(86, 13).Value = -249999877  # M86: -71428587 -> -249999877
$ws.Cells.Item(86, 14).Value = -6199.5715  # N86: -6245.6667 -> -6199.5715
# Row 89
$ws.Cells.Item(89, 8).Value = 55558852  # H89: 50002000 -> 55558852
$ws.Cells.Item(89, 9).Value = 250001000  # I89: 71429710 -> 250001000
$ws.Cells.Item(89, 10).Value = 3953.5715  # J89: 3999.6667 -> 3953.5715
$ws.Cells.Item(89, 11).Value = 1250005000  # K89: 357148550 -> 1250005000
$ws.Cells.Item(89, 12).Value = 19767.8575  # L89: 19998.3335 -> 19767.8575
$ws.Cells.Item(89, 13).Value = -1249999384  # M89: -357142934 -> -1249999384
$ws.Cells.Item(89, 14).Value = -30999.8575  # N89: -31230.3335 -> -30999.8575
# Row 134
$ws.Cells.Item(134, 8).Value = 3198.3333  # H134: 2882.4644 -> 3198.3333
$ws.Cells.Item(134, 9).Value = 1571.0714  # I134: 1441.3334 -> 1571.0714
$ws.Cells.Item(134, 11).Value = 4713.2142  # K134: 4324.0002 -> 4713.2142
$ws.Cells.Item(134, 13).Value = -2178.2142  # M134: -1789.0002 -> -2178.2142
# Row 135
$ws.Cells.Item(135, 8).Value = 43990  # H135: 0 -> 43990
$ws.Cells.Item(135, 10).Value = 43990  # J135: 0 -> 43990
$ws.Cells.Item(135, 12).Value = 43990  # L135: 0 -> 43990
$ws.Cells.Item(135, 14).Value = -54130  # N135: None -> -54130

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Cells.Item(3, 8).Value = 6762.222  # H3: 6543.3335 -> 6762.222
$ws.Cells.Item(3, 9).Value = 2030  # I3: 1696.6666 -> 2030
$ws.Cells.Item(3, 10).Value = 8114.2856  # J3: 8966.666999999999 -> 8114.2856
$ws.Cells.Item(3, 11).Value = 6090  # K3: 5089.9998 -> 6090
$ws.Cells.Item(3, 12).Value = 24342.8568  # L3: 26900.001 -> 24342.8568
$ws.Cells.Item(3, 13).Value = -5978  # M3: -4977.9998 -> -5978
$ws.Cells.Item(3, 14).Value = -24566.8568  # N3: -27124.001 -> -24566.8568
# Row 57
$ws.Cells.Item(57, 8).Value = 3600  # H57: 0 -> 3600
$ws.Cells.Item(57, 9).Value = 3600  # I57: 0 -> 3600
$ws.Cells.Item(57, 11).Value = 10800  # K57: 0 -> 10800
$ws.Cells.Item(57, 13).Value = -10241  # M57: None -> -10241
# Row 76
$ws.Cells.Item(76, 8).Value = 0  # H76: 4000 -> 0
$ws.Cells.Item(76, 9).Value = 0  # I76: 3000 -> 0
$ws.Cells.Item(76, 10).Value = 0  # J76: 5000 -> 0
$ws.Cells.Item(76, 11).Value = 0  # K76: 9000 -> 0
$ws.Cells.Item(76, 12).Value = 0  # L76: 15000 -> 0
$ws.Cells.Item(76, 13).ClearContents()  # M76: -8617 -> (cleared)
$ws.Cells.Item(76, 14).ClearContents()  # N76: -15766 -> (cleared)
# Row 79
$ws.Cells.Item(79, 8).Value = 0  # H79: 4000 -> 0
$ws.Cells.Item(79, 9).Value = 0  # I79: 3000 -> 0
$ws.Cells.Item(79, 10).Value = 0  # J79: 5000 -> 0
$ws.Cells.Item(79, 11).Value = 0  # K79: 9000 -> 0
$ws.Cells.Item(79, 12).Value = 0  # L79: 15000 -> 0
$ws.Cells.Item(79, 13).ClearContents()  # M79: -7674 -> (cleared)
$ws.Cells.Item(79, 14).ClearContents()  # N79: -17652 -> (cleared)
# Row 88
$ws.Cells.Item(88, 8).Value = 0  # H88: 3254 -> 0
$ws.Cells.Item(88, 10).Value = 0  # J88: 3254 -> 0
$ws.Cells.Item(88, 12).Value = 0  # L88: 9762 -> 0
$ws.Cells.Item(88, 14).ClearContents()  # N88: -10618 -> (cleared)
# Row 91
$ws.Cells.Item(91, 8).Value = 0  # H91: 3254 -> 0
$ws.Cells.Item(91, 10).Value = 0  # J91: 3254 -> 0
$ws.Cells.Item(91, 12).Value = 0  # L91: 9762 -> 0
$ws.Cells.Item(91, 14).ClearContents()  # N91: -12726 -> (cleared)
# Row 92
$ws.Cells.Item(92, 8).Value = 2000  # H92: 2126.5715 -> 2000
$ws.Cells.Item(92, 10).Value = 0  # J92: 2221.5 -> 0
$ws.Cells.Item(92, 12).Value = 0  # L92: 6664.5 -> 0
$ws.Cells.Item(92, 14).ClearContents()  # N92: -9160.5 -> (cleared)
# Row 94
$ws.Cells.Item(94, 8).Value = 1024  # H94: 2005.75 -> 1024
$ws.Cells.Item(94, 9).Value = 1024  # I94: 1007.6667 -> 1024
$ws.Cells.Item(94, 10).Value = 0  # J94: 5000 -> 0
$ws.Cells.Item(94, 11).Value = 3072  # K94: 3023.0001 -> 3072
$ws.Cells.Item(94, 12).Value = 0  # L94: 15000 -> 0
$ws.Cells.Item(94, 13).Value = -2396  # M94: -2347.0001 -> -2396
$ws.Cells.Item(94, 14).ClearContents()  # N94: -16352 -> (cleared)
# Row 97
$ws.Cells.Item(97, 8).Value = 830.6667  # H97: 737.1667 -> 830.6667
$ws.Cells.Item(97, 9).Value = 830.6667  # I97: 726.5714 -> 830.6667
$ws.Cells.Item(97, 10).Value = 0  # J97: 741.5294 -> 0
$ws.Cells.Item(97, 11).Value = 2492.0001  # K97: 2179.7142 -> 2492.0001
$ws.Cells.Item(97, 12).Value = 0  # L97: 2224.5882 -> 0
$ws.Cells.Item(97, 13).Value = -1996.0001  # M97: -1683.7142 -> -1996.0001
$ws.Cells.Item(97, 14).ClearContents()  # N97: -3216.5882 -> (cleared)
# Row 100
$ws.Cells.Item(100, 8).Value = 5000  # H100: 5007 -> 5000
$ws.Cells.Item(100, 10).Value = 5000  # J100: 5007 -> 5000
$ws.Cells.Item(100, 12).Value = 15000  # L100: 15021 -> 15000
$ws.Cells.Item(100, 14).Value = -16622  # N100: -16643 -> -16622
# Row 107
$ws.Cells.Item(107, 8).Value = 392.02274  # H107: 392.22726 -> 392.02274
$ws.Cells.Item(107, 9).Value = 421  # I107: 431.91306 -> 421
$ws.Cells.Item(107, 10).Value = 357.25  # J107: 348.7619 -> 357.25
$ws.Cells.Item(107, 11).Value = 1263  # K107: 1295.73918 -> 1263
$ws.Cells.Item(107, 12).Value = 1071.75  # L107: 1046.2857 -> 1071.75
$ws.Cells.Item(107, 13).Value = 657  # M107: 624.26082 -> 657
$ws.Cells.Item(107, 14).Value = -4911.75  # N107: -4886.2857 -> -4911.75
# Row 121
$ws.Cells.Item(121, 8).Value = 791.6667  # H121: 1503.75 -> 791.6667
$ws.Cells.Item(121, 9).Value = 331.25  # I121: 365 -> 331.25
$ws.Cells.Item(121, 10).Value = 1160  # J121: 1883.3334 -> 1160
$ws.Cells.Item(121, 11).Value = 993.75  # K121: 1095 -> 993.75
$ws.Cells.Item(121, 12).Value = 3480  # L121: 5650.0002 -> 3480
$ws.Cells.Item(121, 13).Value = 316.25  # M121: 215 -> 316.25
$ws.Cells.Item(121, 14).Value = -6100  # N121: -8270.0002 -> -6100
# Row 133
$ws.Cells.Item(133, 8).Value = 8051.4287  # H133: 8622.857 -> 8051.4287
$ws.Cells.Item(133, 9).Value = 5090  # I133: 4120 -> 5090
$ws.Cells.Item(133, 11).Value = 15270  # K133: 12360 -> 15270
$ws.Cells.Item(133, 13).Value = -10210  # M133: -7300 -> -10210

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Cells.Item(43, 8).Value = 4533.778  # H43: 6050.6665 -> 4533.778
$ws.Cells.Item(43, 9).Value = 2800.6667  # I43: 3260.8 -> 2800.6667
$ws.Cells.Item(43, 10).Value = 8000  # J43: 20000 -> 8000
$ws.Cells.Item(43, 11).Value = 2800.6667  # K43: 3260.8 -> 2800.6667
$ws.Cells.Item(43, 12).Value = 8000  # L43: 20000 -> 8000
$ws.Cells.Item(43, 13).Value = -2649.6667  # M43: -3109.8 -> -2649.6667
$ws.Cells.Item(43, 14).Value = -8302  # N43: -20302 -> -8302
# Row 102
$ws.Cells.Item(102, 8).Value = 1858.6389  # H102: 1954.8572 -> 1858.6389
$ws.Cells.Item(102, 9).Value = 1710.5714  # I102: 1781.55 -> 1710.5714
$ws.Cells.Item(102, 10).Value = 2065.9333  # J102: 2185.9333 -> 2065.9333
$ws.Cells.Item(102, 11).Value = 1710.5714  # K102: 1781.55 -> 1710.5714
$ws.Cells.Item(102, 12).Value = 2065.9333  # L102: 2185.9333 -> 2065.9333
$ws.Cells.Item(102, 13).Value = -88.57140000000004  # M102: -159.55 -> -88.57140000000004
$ws.Cells.Item(102, 14).Value = -5309.933300000001  # N102: -5429.933300000001 -> -5309.933300000001
# Row 134
$ws.Cells.Item(134, 8).Value = 19730.4  # H134: 20002.2 -> 19730.4
$ws.Cells.Item(134, 10).Value = 19730.4  # J134: 20002.2 -> 19730.4
$ws.Cells.Item(134, 12).Value = 59191.2  # L134: 60006.60000000001 -> 59191.2
$ws.Cells.Item(134, 14).Value = -64261.2  # N134: -65076.60000000001 -> -64261.2

$ws = $wb.Worksheets.Item("LTW")
# Row 101
$ws.Cells.Item(101, 8).Value = 28195  # H101: 29900 -> 28195
$ws.Cells.Item(101, 10).Value = 28195  # J101: 29900 -> 28195
$ws.Cells.Item(101, 12).Value = 28195  # L101: 29900 -> 28195
$ws.Cells.Item(101, 14).Value = -34685  # N101: -36390 -> -34685
# Row 106
$ws.Cells.Item(106, 8).Value = 18287.143  # H106: 19790 -> 18287.143
$ws.Cells.Item(106, 10).Value = 18287.143  # J106: 19790 -> 18287.143
$ws.Cells.Item(106, 12).Value = 18287.143  # L106: 19790 -> 18287.143
$ws.Cells.Item(106, 14).Value = -20811.143  # N106: -22314 -> -20811.143

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Cells.Item(4, 8).Value = 4993.3335  # H4: 0 -> 4993.3335
$ws.Cells.Item(4, 9).Value = 4980  # I4: 0 -> 4980
$ws.Cells.Item(4, 10).Value = 5000  # J4: 0 -> 5000
$ws.Cells.Item(4, 11).Value = 4980  # K4: 0 -> 4980
$ws.Cells.Item(4, 12).Value = 5000  # L4: 0 -> 5000
$ws.Cells.Item(4, 13).Value = -4867  # M4: None -> -4867
$ws.Cells.Item(4, 14).Value = -5226  # N4: None -> -5226
# Row 46
$ws.Cells.Item(46, 8).Value = 196809.67  # H46: 262582.25 -> 196809.67
$ws.Cells.Item(46, 10).Value = 196809.67  # J46: 262582.25 -> 196809.67
$ws.Cells.Item(46, 12).Value = 196809.67  # L46: 262582.25 -> 196809.67
$ws.Cells.Item(46, 14).Value = -197271.67  # N46: -263044.25 -> -197271.67
# Row 97
$ws.Cells.Item(97, 8).Value = 29260  # H97: 29293.334 -> 29260
$ws.Cells.Item(97, 10).Value = 29260  # J97: 29293.334 -> 29260
$ws.Cells.Item(97, 12).Value = 29260  # L97: 29293.334 -> 29260
$ws.Cells.Item(97, 14).Value = -31242  # N97: -31275.334 -> -31242
# Row 103
$ws.Cells.Item(103, 8).Value = 355167.34  # H103: 355200.66 -> 355167.34
$ws.Cells.Item(103, 10).Value = 355167.34  # J103: 355200.66 -> 355167.34
$ws.Cells.Item(103, 12).Value = 355167.34  # L103: 355200.66 -> 355167.34
$ws.Cells.Item(103, 14).Value = -357511.34  # N103: -357544.66 -> -357511.34
# Row 126
$ws.Cells.Item(126, 8).Value = 63707.25  # H126: 72750.92999999999 -> 63707.25
$ws.Cells.Item(126, 9).Value = 125676.5  # I126: 167434.83 -> 125676.5
$ws.Cells.Item(126, 11).Value = 377029.5  # K126: 502304.49 -> 377029.5
$ws.Cells.Item(126, 13).Value = -374559.5  # M126: -499834.49 -> -374559.5
# Row 134
$ws.Cells.Item(134, 8).Value = 196809.67  # H134: 262582.25 -> 196809.67
$ws.Cells.Item(134, 10).Value = 196809.67  # J134: 262582.25 -> 196809.67
$ws.Cells.Item(134, 12).Value = 590429.01  # L134: 787746.75 -> 590429.01
$ws.Cells.Item(134, 14).Value = -595499.01  # N134: -792816.75 -> -595499.01
# Row 136
$ws.Cells.Item(136, 8).Value = 1433.25  # H136: 1541.96 -> 1433.25
$ws.Cells.Item(136, 9).Value = 840.4783  # I136: 916.619 -> 840.4783
$ws.Cells.Item(136, 10).Value = 4160  # J136: 4825 -> 4160
$ws.Cells.Item(136, 11).Value = 2521.4349  # K136: 2749.857 -> 2521.4349
$ws.Cells.Item(136, 12).Value = 12480  # L136: 14475 -> 12480
$ws.Cells.Item(136, 13).Value = 28.5650999999998  # M136: -199.857 -> 28.5650999999998
$ws.Cells.Item(136, 14).Value = -17580  # N136: -19575 -> -17580
